# Update the table style applied to the three summary tables (slides 14-16)
# from the custom "Table_0" style to the built-in style referenced by
# {2F041530-37D1-4869-B407-A6C42A90D82F}.
#
# Table styles cannot be changed by assigning Table.Style/StyleId directly;
# PowerPoint's object model requires the dedicated ApplyStyle method.

$p = $ppt.ActivePresentation

$newStyleId = "{2F041530-37D1-4869-B407-A6C42A90D82F}"

foreach ($slideIndex in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIndex)
    $shape = $slide.Shapes.Item(1)
    if ($shape.HasTable) {
        $shape.Table.ApplyStyle($newStyleId)
    }
}
